# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "50.968.79"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").Value = "2.951.83"

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.NumberFormat = "General"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "380.21"
$c.NumberFormat = "General"
$ws.Range("E5").Value = "  +0.84%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "101.83"
$c.NumberFormat = "General"
$ws.Range("E6").Value = "  +0.24%  "

# Row 7
$ws.Range("E7").Value = "  +1.74%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.585"
$c.NumberFormat = "General"
$ws.Range("E9").Value = "  +0.95%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.24"
$c.NumberFormat = "General"
$ws.Range("E10").Value = "  +0.03%  "

# Row 11
$ws.Range("E11").Value = "  -0.31%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0849"
$c.NumberFormat = "General"
$ws.Range("E12").Value = "  +2.13%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "18.40"
$c.NumberFormat = "General"
$ws.Range("E13").Value = "  +3.53%  "

# Row 14
$ws.Range("D14").Value = "3.411.50"
$ws.Range("E14").Value = "  +0.83%  "

# Row 15
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "12.35"
$c.NumberFormat = "General"
$ws.Range("E15").Value = "  +74.87%  "

# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.75"
$c.NumberFormat = "General"
$ws.Range("E16").Value = "  +6.24%  "

# Row 17
$ws.Range("D17").Value = "2.948.68"
$ws.Range("E17").Value = "  +0.82%  "

# Row 18
$ws.Range("E18").Value = "  +4.54%  "

# Row 19
$ws.Range("D19").Value = "50.984.51"
$ws.Range("E19").Value = "  +0.32%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.07"
$c.NumberFormat = "General"
$ws.Range("E20").Value = "  -2.24%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.34"
$c.NumberFormat = "General"
$ws.Range("E21").Value = "  -0.62%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0953"
$ws.Range("E22").Value = "  +1.05%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.40"
$c.NumberFormat = "General"
$ws.Range("E23").Value = "  +18.94%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "69.57"
$c.NumberFormat = "General"
$ws.Range("E24").Value = "  +2.41%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "266.42"
$c.NumberFormat = "General"
$ws.Range("E25").Value = "  +2.20%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "8.00"
$c.NumberFormat = "General"
$ws.Range("E26").Value = "  -1.29%  "

# Row 27
$ws.Range("E27").Value = "  -0.12%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "25.76"
$c.NumberFormat = "General"
$ws.Range("E28").Value = "  +1.34%  "

# Row 29
$ws.Range("E29").Value = "  -0.79%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.91"
$c.NumberFormat = "General"
$ws.Range("E30").Value = "  -7.95%  "

# Row 31
$ws.Range("E31").Value = "  -5.07%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "10.45"
$c.NumberFormat = "General"
$ws.Range("E32").Value = "  +7.37%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "50.70"
$c.NumberFormat = "General"
$ws.Range("E33").Value = "  +0.02%  "

# Row 34
$ws.Range("E34").Value = "  +0.94%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "33.89"
$c.NumberFormat = "General"
$ws.Range("E35").Value = "  +0.74%  "

# Row 36
$ws.Range("E36").Value = "  -3.25%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.NumberFormat = "General"
$ws.Range("E38").Value = "  +8.13%  "

# Row 39
$ws.Range("E39").Value = "  +2.16%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "16.62"
$c.NumberFormat = "General"
$ws.Range("E40").Value = "  +2.80%  "

# Row 41
$ws.Range("E41").Value = "  +2.89%  "

# Row 42
$ws.Range("E42").Value = "  -4.18%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "118.86"
$c.NumberFormat = "General"
$ws.Range("E43").Value = "  -1.23%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.NumberFormat = "General"
$ws.Range("E44").Value = "  +11.50%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "21.34"
$c.NumberFormat = "General"
$ws.Range("E45").Value = "  +1.61%  "

# Row 46
$ws.Range("E46").Value = "  -0.92%  "

# Row 47
$ws.Range("D47").Value = "2.023.14"
$ws.Range("E47").Value = "  +1.41%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.33"
$c.NumberFormat = "General"
$ws.Range("E48").Value = "  -1.91%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.259"
$c.NumberFormat = "General"
$ws.Range("E49").Value = "  -4.08%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0318"
$c.NumberFormat = "General"
$ws.Range("E50").Value = "  -7.35%  "

# Row 51
$ws.Range("E51").Value = "  +7.13%  "

